# "Missing figures in simulation and the whole results and conclusion"
#
# The "Frec" column (column A) on Sheet1 was actually recorded/typed in Hz,
# but had been labelled/entered as if it were in KHz. Fix the unit: relabel
# the header from "Frec (KHz)" to "Frec (Hz)" and multiply every recorded
# frequency (rows 2-19) by 1000 so the figures are correct.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header: "Frec (KHz)" -> "Frec (Hz)"
$ws.Range("A1").Value = "Frec (Hz)"

# Recorded frequencies were in KHz; convert them to Hz (x1000)
$ws.Range("A2").Value = 1000
$ws.Range("A3").Value = 8000
$ws.Range("A4").Value = 15000
$ws.Range("A5").Value = 20000
$ws.Range("A6").Value = 25000
$ws.Range("A7").Value = 30000
$ws.Range("A8").Value = 35000
$ws.Range("A9").Value = 40000
$ws.Range("A10").Value = 45000
$ws.Range("A11").Value = 50000
$ws.Range("A12").Value = 60000
$ws.Range("A13").Value = 75000
$ws.Range("A14").Value = 85000
$ws.Range("A15").Value = 100000
$ws.Range("A16").Value = 130000
$ws.Range("A17").Value = 200000
$ws.Range("A18").Value = 250000
$ws.Range("A19").Value = 500000

# Leave the cursor on row 2 (whole row selected), matching where the author
# ended up after reviewing/fixing the results
$ws.Rows(2).Select()
